$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Revert "User data 3.0": remove the "budget-type" column (B), shifting
# the 2013/2014/2015/2016 value columns (old C:F) left into B:E.
$ws.Columns.Item(2).Delete()
